# Daily attendance processing - 2025-10-27 23:19:40
# For every "Recorded By" (column G) cell whose value is a comma-separated
# list of recorders, rotate the list so the last entry moves to the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2

    if ($v -like "*,*") {
        $parts = $v -split ", "
        $count = $parts.Count
        $last = $parts[$count - 1]
        $rest = $parts[0..($count - 2)]
        $newParts = @($last) + $rest
        $newVal = $newParts -join ", "
        $cell.Value = $newVal
    }
}
